$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet tracked rating dates in columns B..E (newest-first: B, then C, D, E).
# This edit adds two more recent weeks (Jun_26 and Jun_27) as new columns, pushing
# the previously existing date columns (Jun_15, Jun_13, Jun_10) three slots to the
# right, and backfills the freshly opened slots with "UN" (no rating yet) for
# every existing broker row. It also appends two new broker rows (Benchmark and
# Evercore ISI) that only have data out through the Jun_17 column.

# Insert 3 new columns before column C. This shifts the old C,D,E (and their
# values/styles) to F,G,H, leaving C,D,E empty.
$ws.Range("C1:E1").EntireColumn.Insert()

# Old B1 ("Jun_17") stays in B1; column E1 is now a fresh blank cell that should
# carry the same "Jun_17" label (mirrors B1) per the new header layout.
$ws.Range("E1").Value = "Jun_17"

# New header labels for the two newly added date columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Backfill the newly inserted C,D,E columns for every existing broker row (2-27)
# with "UN" (unchanged / no new rating that week).
$ws.Range("C2:E27").Value = "UN"

# Add the two new broker rows at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
